# mode_conso_global.xlsx — "Encore des chgmts de titres"
#
# The "Year of Treatment" column (B) is dropped entirely, shifting every
# column C:I one slot left (to B:H), and each remaining measure header
# gets a ".global" suffix appended (Inject -> Inject.global, etc.),
# including the former "Total" column which becomes "Total.global".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column B ("Year of Treatment") — shifts C:I left into B:H and
# updates the sheet dimension from A1:I34 to A1:H34 automatically.
$ws.Columns("B:B").Delete()

# Re-label the shifted header row with the ".global" suffix.
$ws.Range("B1").Value = "Inject.global"
$ws.Range("C1").Value = "Smoke / inhale.global"
$ws.Range("D1").Value = "Eat / drink.global"
$ws.Range("E1").Value = "Sniff.global"
$ws.Range("F1").Value = "Other.global"
$ws.Range("G1").Value = "Not known / missing.global"
$ws.Range("H1").Value = "Total.global"
